# Add sensitivity test runs (2050_TM161_SENS_3A_01 and P07-P10) to the ModelRuns sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModelRuns")

# ---------------------------------------------------------------------------
# 1. Prepare the 5 new rows (283-287) by copying the formatting of row 233,
#    which is the analogous "2035" sensitivity-test row and already carries
#    the exact per-column styles used by the new rows.
# ---------------------------------------------------------------------------
$ws.Range("A233:Y233").Copy()
$ws.Range("A283:Y287").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Row 283 : 2050_TM161_SENS_3A_01
#    (values typed in the order: B, M, G, F, then the rest)
# ---------------------------------------------------------------------------
$ws.Cells.Item(283, 2).Value = "2050_TM161_SENS_3A_01"
$ws.Cells.Item(283, 13).Value = "https://app.asana.com/1/11860278793487/task/1211096749996011?focus=true"
$ws.Cells.Item(283, 7).Value = "SENS_v3A\2050"
$ws.Cells.Item(283, 6).Value = "2050 NoProject_network"

$ws.Cells.Item(283, 1).Value = 2050
$ws.Cells.Item(283, 3).Value = "RTP_2025Sens"
$ws.Cells.Item(283, 4).Value = "FBP"
$ws.Cells.Item(283, 5).Value = "Integrated Model - highway widening"
$ws.Cells.Item(283, 8).Value = "PBA50Plus_sensitivity_V3A"
$ws.Cells.Item(283, 9).Value = "sensitivity_longRunInducedDemand"
$ws.Cells.Item(283, 11).Value = "BlueprintNetworks_v35\net_2050_Baseline"
$ws.Cells.Item(283, 12).Value = "model3-b"
$ws.Cells.Item(283, 14).Value = 16.47
$ws.Cells.Item(283, 15).Value = "na"
$ws.Cells.Item(283, 16).Value = "na"
$ws.Cells.Item(283, 20).Value = -0.455
$ws.Cells.Item(283, 21).Value = 5
$ws.Cells.Item(283, 22).Value = 55
$ws.Cells.Item(283, 23).Value = 0
$ws.Cells.Item(283, 24).Value = 108
$ws.Cells.Item(283, 25).Value = "NoProject_network"

# ---------------------------------------------------------------------------
# 3. Update the existing 2035 analog row (226) description to match the new
#    naming convention.
# ---------------------------------------------------------------------------
$ws.Cells.Item(226, 6).Value = "2035 NoProject_network"

# ---------------------------------------------------------------------------
# 4. Fill the "directory" (B) column for rows 284-287 first (top to bottom),
#    matching the order the shared strings were originally recorded in.
# ---------------------------------------------------------------------------
$ws.Cells.Item(284, 2).Value = "2050_TM161_SENS_3A_P07"
$ws.Cells.Item(285, 2).Value = "2050_TM161_SENS_3A_P08"
$ws.Cells.Item(286, 2).Value = "2050_TM161_SENS_3A_P09"
$ws.Cells.Item(287, 2).Value = "2050_TM161_SENS_3A_P10"

# ---------------------------------------------------------------------------
# 5. Row 284 : 2050_TM161_SENS_3A_P07 (remaining columns; K then M)
# ---------------------------------------------------------------------------
$ws.Cells.Item(284, 11).Value = "BlueprintNetworks_v37_for_SENS\2050_v7_eastCC_Oakland_34LaneAdd"
$ws.Cells.Item(284, 13).Value = "https://app.asana.com/1/11860278793487/project/1203667963226602/task/1211118174941124?focus=true"

$ws.Cells.Item(284, 1).Value = 2050
$ws.Cells.Item(284, 3).Value = "RTP_2025Sens"
$ws.Cells.Item(284, 4).Value = "FBP"
$ws.Cells.Item(284, 5).Value = "Integrated Model - highway widening"
$ws.Cells.Item(284, 6).Value = "eastCC_Oakland_34LaneAdd"
$ws.Cells.Item(284, 7).Value = "SENS_v3A\2050"
$ws.Cells.Item(284, 8).Value = "PBA50Plus_sensitivity_V3A"
$ws.Cells.Item(284, 9).Value = "sensitivity_longRunInducedDemand"
$ws.Cells.Item(284, 12).Value = "model3-b"
$ws.Cells.Item(284, 14).Value = 16.47
$ws.Cells.Item(284, 15).Value = "na"
$ws.Cells.Item(284, 16).Value = "na"
$ws.Cells.Item(284, 20).Value = -0.455
$ws.Cells.Item(284, 21).Value = 5
$ws.Cells.Item(284, 22).Value = 55
$ws.Cells.Item(284, 23).Value = 0
$ws.Cells.Item(284, 24).Value = 108
$ws.Cells.Item(284, 25).Value = "eastCC_Oakland_34LaneAdd"

# ---------------------------------------------------------------------------
# 6. Row 285 : 2050_TM161_SENS_3A_P08 (remaining columns; M then K)
# ---------------------------------------------------------------------------
$ws.Cells.Item(285, 13).Value = "https://app.asana.com/1/11860278793487/project/1203667963226602/task/1211118174941128?focus=true"
$ws.Cells.Item(285, 11).Value = "BlueprintNetworks_v37_for_SENS\2050_v8_eastALACC_Oakland_34LaneAdd"

$ws.Cells.Item(285, 1).Value = 2050
$ws.Cells.Item(285, 3).Value = "RTP_2025Sens"
$ws.Cells.Item(285, 4).Value = "FBP"
$ws.Cells.Item(285, 5).Value = "Integrated Model - highway widening"
$ws.Cells.Item(285, 6).Value = "eastALACC_Oakland_34LaneAdd"
$ws.Cells.Item(285, 7).Value = "SENS_v3A\2050"
$ws.Cells.Item(285, 8).Value = "PBA50Plus_sensitivity_V3A"
$ws.Cells.Item(285, 9).Value = "sensitivity_longRunInducedDemand"
$ws.Cells.Item(285, 12).Value = "model3-c"
$ws.Cells.Item(285, 14).Value = 16.47
$ws.Cells.Item(285, 15).Value = "na"
$ws.Cells.Item(285, 16).Value = "na"
$ws.Cells.Item(285, 20).Value = -0.455
$ws.Cells.Item(285, 21).Value = 5
$ws.Cells.Item(285, 22).Value = 55
$ws.Cells.Item(285, 23).Value = 0
$ws.Cells.Item(285, 24).Value = 108
$ws.Cells.Item(285, 25).Value = "eastALACC_Oakland_34LaneAdd"

# ---------------------------------------------------------------------------
# 7. Row 286 : 2050_TM161_SENS_3A_P09 (remaining columns; M then K)
# ---------------------------------------------------------------------------
$ws.Cells.Item(286, 13).Value = "https://app.asana.com/1/11860278793487/project/1203667963226602/task/1211118598663388?focus=true"
$ws.Cells.Item(286, 11).Value = "BlueprintNetworks_v37_for_SENS\2050_v9_ALA_SMSCL_34LaneAdd"

$ws.Cells.Item(286, 1).Value = 2050
$ws.Cells.Item(286, 3).Value = "RTP_2025Sens"
$ws.Cells.Item(286, 4).Value = "FBP"
$ws.Cells.Item(286, 5).Value = "Integrated Model - highway widening"
$ws.Cells.Item(286, 6).Value = "ALA_SMSCL_34LaneAdd"
$ws.Cells.Item(286, 7).Value = "SENS_v3A\2050"
$ws.Cells.Item(286, 8).Value = "PBA50Plus_sensitivity_V3A"
$ws.Cells.Item(286, 9).Value = "sensitivity_longRunInducedDemand"
$ws.Cells.Item(286, 12).Value = "model2-c"
$ws.Cells.Item(286, 14).Value = 16.47
$ws.Cells.Item(286, 15).Value = "na"
$ws.Cells.Item(286, 16).Value = "na"
$ws.Cells.Item(286, 20).Value = -0.455
$ws.Cells.Item(286, 21).Value = 5
$ws.Cells.Item(286, 22).Value = 55
$ws.Cells.Item(286, 23).Value = 0
$ws.Cells.Item(286, 24).Value = 108
$ws.Cells.Item(286, 25).Value = "ALA_SMSCL_34LaneAdd"

# ---------------------------------------------------------------------------
# 8. Row 287 : 2050_TM161_SENS_3A_P10 (remaining columns; M then K)
# ---------------------------------------------------------------------------
$ws.Cells.Item(287, 13).Value = "https://app.asana.com/1/11860278793487/project/1203667963226602/task/1211118602655911?focus=true"
$ws.Cells.Item(287, 11).Value = "BlueprintNetworks_v37_for_SENS\2050_v10_Vallejo_SF_34LaneAdd"

$ws.Cells.Item(287, 1).Value = 2050
$ws.Cells.Item(287, 3).Value = "RTP_2025Sens"
$ws.Cells.Item(287, 4).Value = "FBP"
$ws.Cells.Item(287, 5).Value = "Integrated Model - highway widening"
$ws.Cells.Item(287, 6).Value = "Vallejo_SF_34LaneAdd"
$ws.Cells.Item(287, 7).Value = "SENS_v3A\2050"
$ws.Cells.Item(287, 8).Value = "PBA50Plus_sensitivity_V3A"
$ws.Cells.Item(287, 9).Value = "sensitivity_longRunInducedDemand"
$ws.Cells.Item(287, 12).Value = "model2-c"
$ws.Cells.Item(287, 14).Value = 16.47
$ws.Cells.Item(287, 15).Value = "na"
$ws.Cells.Item(287, 16).Value = "na"
$ws.Cells.Item(287, 20).Value = -0.455
$ws.Cells.Item(287, 21).Value = 5
$ws.Cells.Item(287, 22).Value = 55
$ws.Cells.Item(287, 23).Value = 0
$ws.Cells.Item(287, 24).Value = 108
$ws.Cells.Item(287, 25).Value = "Vallejo_SF_34LaneAdd"

# ---------------------------------------------------------------------------
# 8. Column width tweaks (status column narrower, network column much wider
#    to accommodate the long new path strings).
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 7.418
$ws.Columns.Item(11).ColumnWidth = 69.585

# ---------------------------------------------------------------------------
# 9. Update the view so the frozen pane / selection reflect scrolling down to
#    the newly added rows.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 270
$ws.Range("M290").Select()
